$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: paragraph "3. [Primary] Remove the the last two logs from
# archive.< as they haven't been shipped yet."
#   -> "3. [Primary] Remove the last two logs from archive. < as they
#       haven't been shipped yet."
# (removes the duplicated "the" and adds a space before the "<")
# -----------------------------------------------------------------------
$find1a = $d.Content.Find
$find1a.Execute("the the last two logs from ", $true, $false, $false, $false, $false, $true, 1, $false, "the last two logs from ", 2) | Out-Null

$find1b = $d.Content.Find
$find1b.Execute("archive.<", $true, $false, $false, $false, $false, $true, 1, $false, "archive. <", 2) | Out-Null

# -----------------------------------------------------------------------
# Change 2: "controlfile" -> "control file" but only in the two step-6
# paragraphs (other occurrences elsewhere in the document stay as-is).
# -----------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Create a new standby controlfile") {
        $f = $p.Range.Find
        $f.Execute("controlfile", $true, $false, $false, $false, $false, $true, 1, $false, "control file", 2) | Out-Null
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "alter database create standby controlfile") {
        $f = $p.Range.Find
        $f.Execute("controlfile", $true, $false, $false, $false, $false, $true, 1, $false, "control file", 2) | Out-Null
        break
    }
}

# -----------------------------------------------------------------------
# Change 3: add two new paragraphs ("Or" / "RESTORE DATABASE FROM BACKUP
# LOCATION '/path/to/backup_files';") right after the
# "RMAN> catalog start with '/home/oracle/raj';" paragraph.
# -----------------------------------------------------------------------
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "catalog start with '/home/oracle/raj'") {
        $anchorPara = $p
        break
    }
}
if ($anchorPara -ne $null) {
    $anchorPara.Range.InsertParagraphAfter()
    $newPara1 = $anchorPara.Next()
    $newPara1.Range.Text = "Or"

    $newPara1.Range.InsertParagraphAfter()
    $newPara2 = $newPara1.Next()
    $newPara2.Range.Text = "RESTORE DATABASE FROM BACKUP LOCATION '/path/to/backup_files';"
}

# -----------------------------------------------------------------------
# Change 4: the three consecutive empty paragraphs (bold formatting mark
# only) that follow "... i.e recovering the database" collapse into a
# single empty paragraph whose paragraph mark simply carries
# <w:spacing w:after="0"/>.
# -----------------------------------------------------------------------
$anchorPara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "recovering the database") {
        $anchorPara2 = $p
        break
    }
}
if ($anchorPara2 -ne $null) {
    $emptyPara = $anchorPara2.Next()
    # Remove two of the three empty paragraphs, leaving one behind.
    $emptyPara.Range.Delete()
    $emptyPara.Range.Delete()
    # Clear the inherited bold paragraph-mark formatting and set the
    # paragraph spacing to match the target.
    $emptyPara.Style = $d.Styles.Item("Normal")
    $emptyPara.SpaceAfter = 0
}
